$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 145, shifting existing rows 145..213 down to 146..214.
$ws.Rows.Item(145).Insert()

# Populate the newly inserted row 145 with the new price record.
$ws.Range("A145").Value = 5
$ws.Range("B145").Value = "Macroferia Regional de Talca"
$ws.Range("C145").Value = "Maule"
$ws.Range("D145").Value = 44553
$ws.Range("E145").Value = 7
$ws.Range("F145").Value = 100112006
$ws.Range("G145").Value = "Repollo"
$ws.Range("H145").Value = "Crespo record"
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 4000
$ws.Range("K145").Value = 700
$ws.Range("L145").Value = 700
$ws.Range("M145").Value = 700
$ws.Range("N145").Value = "`$/unidad"
$ws.Range("O145").Value = "Región del Maule"
$ws.Range("P145").Value = 700
$ws.Range("Q145").Value = 1
$ws.Range("R145").Value = "Hortaliza"
